# SCCB-S13276: fix the footer copyright line.
# "Copyright © 2011, Inc All Rights Reserved." -> "Copyright © 2011, All Rights Reserved."
# (drops the stray "Inc" that was spanning its own run/proofErr pair)

$d = $word.ActiveDocument

$copyrightSign = [char]0x00A9

$oldText = "Copyright " + $copyrightSign + " 2011, Inc All Rights Reserved."
$newText = "Copyright " + $copyrightSign + " 2011, All Rights Reserved."

$replaced = $false

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections.Item($s)
    for ($i = 1; $i -le 3; $i++) {
        $footer = $section.Footers.Item($i)
        if ($footer.Exists) {
            $rng = $footer.Range.Duplicate
            $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                                        $true, 1, $false, $newText, 2)
            if ($found) {
                $replaced = $true
            }
        }
    }
}

Write-Output ("Replaced copyright text: " + $replaced)
